$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Start period", "End period" and "Apartment" columns (old D:F) -
# the revenue/expense tracker is moving to an invoice-driven flow instead of
# period tracking.
$ws.Range("D1:F1").EntireColumn.Delete()

# Rename the "Revenue type" header (old B, now still B) to "Invoice".
$ws.Range("B1").Value = "Invoice"

# Land the selection on A2 like the author's updated template.
$ws.Range("A2").Select()
